$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns P and Q: copy formatting from the existing last header/data-col style, then set header values
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# Data rows 2-25: columns B..I (2..9) get updated values; P (16) and Q (17) are new zero-valued columns
$rows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25)
$data = @(
    @(3.469011592393883,1.021860864717723,0.345883679821597,1.353683879945791,6.452836240317794,0,0.01113859335798506,0.00375459660703914),
    @(3.007426218724959,0.882495996475285,0.3031011142505378,1.164929209104798,5.639920025744971,0,0.006523286070125867,0.001126947876521012),
    @(2.726118036786715,0.7987546882678487,0.2768169867672725,1.050780429963453,5.143305393537901,0,0.004249831763293743,0.0004789386768653259),
    @(2.611839370048244,0.7665372401730508,0.264875625703155,1.004570121369113,4.927549603468947,0,0.003445900418873382,0.0005003636661573552),
    @(2.592798102614665,0.7630585986810274,0.2614141850518621,0.9968384907395915,4.874935878556641,0,0.003313688035798901,0.0006041873337165882),
    @(2.724311823131302,0.8033689126937134,0.2725864070509942,1.049928209706778,5.094198923027335,0,0.004219396111079687,0.0006936530302574795),
    @(3.308902216886906,0.9800930657349056,0.3255981066706113,1.287835512273062,6.109457475727254,0,0.009384658035615973,0.002848585562690076),
    @(4.477765696643075,1.335064853993742,0.4377660521217308,1.77217049457434,8.218493594963434,0,0.02468351707386152,0.01563464254944069),
    @(5.353061617035053,1.609228068759649,0.4974241334507497,2.036752382280596,9.532796148470084,0,0.03876215825391593,0.03143500498223517),
    @(5.717335597828821,1.709277877030047,0.3445206894515849,1.332578271195118,7.986060767815474,0,0.05310200936765241,0.03352720109775476),
    @(5.842270279576212,1.730270927171375,0.2347933616326117,0.8203140205408346,6.621554741140756,0,0.08729021839531015,0.03190847307126266),
    @(5.789015640977368,1.700020546434587,0.1459429037387778,0.4279166647295156,5.255122895374114,0,0.1380835277082326,0.02770712639422612),
    @(5.672343455190401,1.658360723366513,0.09624954436909405,0.2273289108355385,4.316174672005275,0,0.1840199874014132,0.02400547098758743),
    @(5.605746308507321,1.639288884221742,0.08503719055227066,0.1870956428553896,4.060381367476481,0,0.1953766633514675,0.02265491307935985),
    @(5.241291209941835,1.533150720907713,0.08236658395773588,0.1765754033928886,3.822440562515709,0,0.178636741573186,0.01743067647423935),
    @(5.023843307058655,1.47365098140466,0.103314271656906,0.259660947893515,4.112457099653994,0,0.1394225125439164,0.01531757462804251),
    @(4.907622235161853,1.443278390973262,0.1599894076873483,0.4968619822184834,4.98317967055624,0,0.08763874446295006,0.01520015362055194),
    @(4.882414802345409,1.451179838986491,0.2554309197415137,0.9419771787818831,6.306661355478184,0,0.04642947035745948,0.01744009710503036),
    @(5.11955239877534,1.551843424259403,0.4671969071152517,1.960633764796171,9.021661065950155,0,0.03452929122076576,0.02703594622133565),
    @(5.807672759244838,1.772013572664036,0.550618505573027,2.335558391679655,10.46718538230971,0,0.0489693828398412,0.04308686461131561),
    @(6.263163369224515,1.911761852741847,0.6008153418076176,2.532010152372891,11.36734425022314,0,0.05902080928144038,0.0551487315617889),
    @(6.020185172135825,1.830722725641294,0.579308618000681,2.427235178077282,10.94508252591754,0,0.05362660949708875,0.04844006406209189),
    @(5.112870350126343,1.542418162858837,0.491092662063437,2.03954418003454,9.28438459272445,0,0.03539875506675561,0.02717145967847934),
    @(4.157600708483187,1.245914403592906,0.3994381108935698,1.638158791497887,7.554878834579938,0,0.01981917451822257,0.0112904188781382)
)

for ($idx = 0; $idx -lt $rows.Length; $idx++) {
    $r = $rows[$idx]
    $vals = $data[$idx]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 2).Value = $vals[$i]
    }
    $ws.Cells.Item($r, 16).Value = 0
    $ws.Cells.Item($r, 17).Value = 0
}
